$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates: rows 2-3 move from 2022-10-04 (44838) to 2022-09-28 (44832)
# and rows 6-7 move from 2022-09-28 (44832) to 2022-10-04 (44838)
$ws.Range("D2").Value = 44832
$ws.Range("D3").Value = 44832
$ws.Range("D6").Value = 44838
$ws.Range("D7").Value = 44838
